$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 2-35 (columns A-D). Columns E/F/G stay "nan"/"nan"/"PUERTOLIBRE".
$data = @(
    @(2, 1233, "172.28.255.41", "MADRI-MADR-H-01-DAAS", "xe-0/0/14"),
    @(3, 1234, "172.28.255.41", "MADRI-MADR-H-01-DAAS", "xe-0/0/15"),
    @(4, 1235, "172.28.255.41", "MADRI-MADR-H-01-DAAS", "xe-0/0/16"),
    @(5, 1236, "172.28.255.41", "MADRI-MADR-H-01-DAAS", "xe-0/0/17"),
    @(6, 1237, "172.28.255.41", "MADRI-MADR-H-01-DAAS", "xe-0/0/18"),
    @(7, 1238, "172.28.255.41", "MADRI-MADR-H-01-DAAS", "xe-0/0/19"),
    @(8, 1240, "172.28.255.41", "MADRI-MADR-H-01-DAAS", "xe-0/0/20"),
    @(9, 1241, "172.28.255.41", "MADRI-MADR-H-01-DAAS", "xe-0/0/21"),
    @(10, 1242, "172.28.255.41", "MADRI-MADR-H-01-DAAS", "xe-0/0/22"),
    @(11, 1243, "172.28.255.41", "MADRI-MADR-H-01-DAAS", "xe-0/0/23"),
    @(12, 1244, "172.28.255.41", "MADRI-MADR-H-01-DAAS", "xe-0/0/24"),
    @(13, 1245, "172.28.255.41", "MADRI-MADR-H-01-DAAS", "xe-0/0/25"),
    @(14, 1246, "172.28.255.41", "MADRI-MADR-H-01-DAAS", "xe-0/0/26"),
    @(15, 1247, "172.28.255.41", "MADRI-MADR-H-01-DAAS", "xe-0/0/27"),
    @(16, 1248, "172.28.255.41", "MADRI-MADR-H-01-DAAS", "xe-0/0/28"),
    @(17, 1249, "172.28.255.41", "MADRI-MADR-H-01-DAAS", "xe-0/0/29"),
    @(18, 1251, "172.28.255.41", "MADRI-MADR-H-01-DAAS", "xe-0/0/30"),
    @(19, 1252, "172.28.255.41", "MADRI-MADR-H-01-DAAS", "xe-0/0/31"),
    @(20, 1253, "172.28.255.41", "MADRI-MADR-H-01-DAAS", "xe-0/0/32"),
    @(21, 1254, "172.28.255.41", "MADRI-MADR-H-01-DAAS", "xe-0/0/33"),
    @(22, 1255, "172.28.255.41", "MADRI-MADR-H-01-DAAS", "xe-0/0/34"),
    @(23, 1257, "172.28.255.41", "MADRI-MADR-H-01-DAAS", "xe-0/0/36"),
    @(24, 1258, "172.28.255.41", "MADRI-MADR-H-01-DAAS", "xe-0/0/37"),
    @(25, 1259, "172.28.255.41", "MADRI-MADR-H-01-DAAS", "xe-0/0/38"),
    @(26, 1260, "172.28.255.41", "MADRI-MADR-H-01-DAAS", "xe-0/0/39"),
    @(27, 1262, "172.28.255.41", "MADRI-MADR-H-01-DAAS", "xe-0/0/40"),
    @(28, 1263, "172.28.255.41", "MADRI-MADR-H-01-DAAS", "xe-0/0/41"),
    @(29, 1264, "172.28.255.41", "MADRI-MADR-H-01-DAAS", "xe-0/0/42"),
    @(30, 1265, "172.28.255.41", "MADRI-MADR-H-01-DAAS", "xe-0/0/43"),
    @(31, 1266, "172.28.255.41", "MADRI-MADR-H-01-DAAS", "xe-0/0/44"),
    @(32, 1267, "172.28.255.41", "MADRI-MADR-H-01-DAAS", "xe-0/0/45"),
    @(33, 1268, "172.28.255.41", "MADRI-MADR-H-01-DAAS", "xe-0/0/46"),
    @(34, 1269, "172.28.255.41", "MADRI-MADR-H-01-DAAS", "xe-0/0/47"),
    @(35, 1270, "172.28.255.41", "MADRI-MADR-H-01-DAAS", "xe-0/0/48")
)

foreach ($r in $data) {
    $rowNum = $r[0]
    $ws.Cells.Item($rowNum, 1).Value = $r[1]
    $ws.Cells.Item($rowNum, 2).Value = $r[2]
    $ws.Cells.Item($rowNum, 3).Value = $r[3]
    $ws.Cells.Item($rowNum, 4).Value = $r[4]
}

# Remove the now-obsolete trailing rows (36-47), shifting nothing below them up
# (there is nothing below) and shrinking the used range to A1:G35.
$ws.Range("A36:G47").EntireRow.Delete() | Out-Null
